$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the data row for LIDIS JOHANA TOVAR PATERNINA (worker no longer in this
# batch of the Estado de Cuenta); remaining rows below shift up by one.
$ws.Rows("17").Delete()

# Update the totals: new "Valor Mora" total and worker count for this update.
$ws.Range("E11").Value = 113880
$ws.Range("C13").Value = 2

# Update the mora period shown for each remaining worker (2507 -> 2508).
$ws.Range("E16").Value = "2508"
$ws.Range("E17").Value = "2508"

# Update the per-worker "Valor Mora" amounts for the new period.
$ws.Range("F16").Value = 56940
$ws.Range("F17").Value = 56940
